$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Yıllar"
$ws.Range("A1").HorizontalAlignment = -4108  # xlCenter
